# Append " (Changed main)" to the end of the first paragraph's text, as
# four separate runs:
#   "This is a Microsoft word document."   (existing run, untouched)
#   " ("                                   (new run)
#   "Changed main"                         (new run)
#   ")"                                    (new run)
#
# A plain Range.InsertAfter() call would just grow the existing run's
# <w:t> (adjacent runs with identical formatting get coalesced), so each
# chunk is instead typed into a fresh paragraph created right after the
# insertion point and then folded back in by deleting the paragraph mark
# that separates it from its predecessor. That merge-by-deleting-the-pilcrow
# step is what keeps the two chunks as distinct <w:r> runs instead of
# re-combining them into one.

$d = $word.ActiveDocument

function Insert-TextAsNewRun($doc, $pos, $text) {
    # Split the paragraph right after $pos into two paragraphs.
    $splitPoint = $doc.Range($pos, $pos)
    $splitPoint.InsertParagraphAfter()

    # Put the new text into the (now separate) paragraph that follows.
    $newParaRange = $doc.Range($pos + 1, $pos + 1)
    $newParaRange.Paragraphs.Item(1).Range.InsertAfter($text)

    # Re-join that paragraph onto the previous one by deleting the
    # paragraph mark at $pos..$pos+1 - this leaves the typed text as its
    # own run rather than merging it back into the preceding run.
    $mark = $doc.Range($pos, $pos + 1)
    $mark.Delete()

    return $pos + $text.Length
}

$firstPara = $d.Paragraphs.Item(1)
$r = $firstPara.Range
# Exclude the trailing paragraph mark from the range so End points right
# after the visible text ("This is a Microsoft word document.").
$r.End = $r.End - 1
$pos = $r.End

$pos = Insert-TextAsNewRun $d $pos " ("
$pos = Insert-TextAsNewRun $d $pos "Changed main"
$pos = Insert-TextAsNewRun $d $pos ")"
